$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new numeric-looking values must stay as text
# (matches source data where these are stored as plain strings, not numbers)
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '57.814.22'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '3.116.33'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '529.98'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D6').Value = '137.93'
$ws.Range('E6').Value = '  -2.21%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E8').Value = '  +3.21%  '
$ws.Range('D9').Value = '7.26'
$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('E10').Value = '  -2.59%  '
$ws.Range('D11').Value = '0.407'
$ws.Range('E11').Value = '  +1.30%  '
$ws.Range('D12').Value = '3.649.25'
$ws.Range('E12').Value = '  -1.23%  '
$ws.Range('E13').Value = '  +1.24%  '
$ws.Range('D14').Value = '25.40'
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('E15').Value = '  -2.52%  '
$ws.Range('D16').Value = '57.796.13'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').Value = '3.106.20'
$ws.Range('E17').Value = '  -1.73%  '
$ws.Range('E18').Value = '  -2.27%  '
$ws.Range('D19').Value = '12.55'
$ws.Range('E19').Value = '  -2.42%  '
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').Value = '350.47'
$ws.Range('E21').Value = '  -1.66%  '
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('D23').Value = '68.91'
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('E24').Value = '  -1.65%  '
$ws.Range('D25').Value = '0.166'
$ws.Range('E25').Value = '  -2.42%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('D27').Value = '0.0₃0867'
$ws.Range('E27').Value = '  -7.97%  '
$ws.Range('D28').Value = '7.19'
$ws.Range('E28').Value = '  -3.86%  '
$ws.Range('E29').Value = '  -2.39%  '
$ws.Range('E30').Value = '  -5.42%  '
$ws.Range('D31').Value = '21.19'
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('D32').Value = '4.93'
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('E33').Value = '  -5.47%  '
$ws.Range('D34').Value = '158.79'
$ws.Range('E34').Value = '  +0.71%  '
$ws.Range('E35').Value = '  -2.91%  '
$ws.Range('D36').Value = '25.76'
$ws.Range('E36').Value = '  -2.39%  '
$ws.Range('E37').Value = '  -3.20%  '
$ws.Range('D38').Value = '1.65'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('D39').Value = '0.0667'
$ws.Range('E39').Value = '  -1.02%  '
$ws.Range('D40').Value = '3.99'
$ws.Range('E40').Value = '  -2.53%  '
$ws.Range('E41').Value = '  -1.65%  '
$ws.Range('E42').Value = '  +1.06%  '
$ws.Range('D43').Value = '2.387.64'
$ws.Range('E43').Value = '  +2.25%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('B45').Value = 'RenzoRestakedETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D45').Value = '3.152.33'
$ws.Range('E45').Value = '  -1.29%  '
$ws.Range('E46').Value = '  -3.43%  '
$ws.Range('D47').Value = '0.958'
$ws.Range('E47').Value = '  -4.73%  '
$ws.Range('D48').Value = '6.02'
$ws.Range('E48').Value = '  -0.99%  '
$ws.Range('D49').Value = '19.66'
$ws.Range('E49').Value = '  -3.52%  '
$ws.Range('D50').Value = '0.734'
$ws.Range('E50').Value = '  -2.96%  '
$ws.Range('D51').Value = '0.0910'
$ws.Range('E51').Value = '  +1.32%  '
